$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 959.7778
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 1091.1428
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 3273.4284
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -3609.4284
$ws.Range("H40").Value = 3879.311
$ws.Range("I40").Value = 2965.9744
$ws.Range("J40").Value = 9816
$ws.Range("K40").Value = 2965.9744
$ws.Range("L40").Value = 9816
$ws.Range("M40").Value = -2790.9744
$ws.Range("N40").Value = -10166
$ws.Range("H62").Value = 7465.273
$ws.Range("I62").Value = 6366.8335
$ws.Range("J62").Value = 8783.4
$ws.Range("K62").Value = 6366.8335
$ws.Range("L62").Value = 8783.4
$ws.Range("M62").Value = -5742.8335
$ws.Range("N62").Value = -10031.4
$ws.Range("H65").Value = 7465.273
$ws.Range("I65").Value = 6366.8335
$ws.Range("J65").Value = 8783.4
$ws.Range("K65").Value = 31834.1675
$ws.Range("L65").Value = 43917
$ws.Range("M65").Value = -28714.1675
$ws.Range("N65").Value = -50157
$ws.Range("H96").Value = 931.6667
$ws.Range("I96").Value = 50.5
$ws.Range("J96").Value = 1372.25
$ws.Range("K96").Value = 151.5
$ws.Range("L96").Value = 4116.75
$ws.Range("M96").Value = 1221.5
$ws.Range("H100").Value = 1599
$ws.Range("I100").Value = 1498.75
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1498.75
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -957.75
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("H127").Value = 1238.6
$ws.Range("I127").Value = 1238.6
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 3715.8
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 1244.2
$ws.Range("H132").Value = 19647.637
$ws.Range("I132").Value = 19647.637
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 58942.91099999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -56412.91099999999
$ws.Range("H135").Value = 586.5333000000001
$ws.Range("I135").Value = 586.5333000000001
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5278.7997
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2743.7997
$ws.Range("H137").Value = 2649.2454
$ws.Range("I137").Value = 1714
$ws.Range("J137").Value = 3549.8518
$ws.Range("K137").Value = 5142
$ws.Range("L137").Value = 10649.5554
$ws.Range("M137").Value = -2592
$ws.Range("H141").Value = 1854.2222
$ws.Range("I141").Value = 1073.5
$ws.Range("J141").Value = 8100
$ws.Range("K141").Value = 3220.5
$ws.Range("L141").Value = 24300
$ws.Range("M141").Value = 1959.5
$ws.Range("N141").Value = -34660
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2290.875
$ws.Range("I45").Value = 1920.7778
$ws.Range("J45").Value = 2766.7144
$ws.Range("K45").Value = 1920.7778
$ws.Range("L45").Value = 2766.7144
$ws.Range("M45").Value = -1543.7778
$ws.Range("N45").Value = -3520.7144
$ws.Range("H61").Value = 4125.7144
$ws.Range("I61").Value = 4125.7144
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4125.7144
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3913.7144
$ws.Range("H95").Value = 46083.2
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 46083.2
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 46083.2
$ws.Range("N95").Value = -51575.2
$ws.Range("H102").Value = 4345.75
$ws.Range("I102").Value = 3040.923
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 3040.923
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -1418.923
$ws.Range("H110").Value = 816
$ws.Range("I110").Value = 804
$ws.Range("J110").Value = 900
$ws.Range("K110").Value = 804
$ws.Range("L110").Value = 900
$ws.Range("M110").Value = 1241
$ws.Range("N110").Value = -4990
$ws.Range("H136").Value = 4125.7144
$ws.Range("I136").Value = 4125.7144
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12377.1432
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9827.143199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3971.25
$ws.Range("I105").Value = 3961.6667
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3961.6667
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -2214.6667
$ws.Range("H107").Value = 4177.467
$ws.Range("I107").Value = 2833
$ws.Range("J107").Value = 7874.75
$ws.Range("K107").Value = 2833
$ws.Range("L107").Value = 7874.75
$ws.Range("M107").Value = -913
$ws.Range("H134").Value = 1234
$ws.Range("I134").Value = 1234
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3702
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1167

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 12271
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 12271
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 12271
$ws.Range("N64").Value = -12767
$ws.Range("H67").Value = 12271
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 12271
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 12271
$ws.Range("N67").Value = -13987
$ws.Range("H81").Value = 49999.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 49999.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 49999.5
$ws.Range("N81").Value = -51995.5
$ws.Range("H84").Value = 49999.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 49999.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 149998.5
$ws.Range("N84").Value = -159982.5
$ws.Range("H107").Value = 287.91666
$ws.Range("I107").Value = 136.125
$ws.Range("J107").Value = 591.5
$ws.Range("K107").Value = 136.125
$ws.Range("L107").Value = 591.5
$ws.Range("M107").Value = 1783.875
$ws.Range("H132").Value = 3829.25
$ws.Range("I132").Value = 1464.4286
$ws.Range("J132").Value = 7140
$ws.Range("K132").Value = 4393.2858
$ws.Range("L132").Value = 21420
$ws.Range("M132").Value = -1863.2858
$ws.Range("N132").Value = -26480
$ws.Range("H134").Value = 1608.2
$ws.Range("I134").Value = 1638.375
$ws.Range("J134").Value = 1487.5
$ws.Range("K134").Value = 4915.125
$ws.Range("L134").Value = 4462.5
$ws.Range("M134").Value = -2380.125
$ws.Range("N134").Value = -9532.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2750
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 7500
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = -7040
$ws.Range("N51").Value = -9920
$ws.Range("H95").Value = 11666.667
$ws.Range("I95").Value = 10000
$ws.Range("J95").Value = 15000
$ws.Range("K95").Value = 30000
$ws.Range("L95").Value = 45000
$ws.Range("M95").Value = -27941
$ws.Range("N95").Value = -49118
$ws.Range("H120").Value = 10457
$ws.Range("I120").Value = 10457
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 31371
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -26533

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 195.05
$ws.Range("I2").Value = 75.72727
$ws.Range("J2").Value = 340.8889
$ws.Range("K2").Value = 75.72727
$ws.Range("L2").Value = 340.8889
$ws.Range("M2").Value = 37.27273
$ws.Range("N2").Value = -566.8888999999999
$ws.Range("H31").Value = 3599.3333
$ws.Range("I31").Value = 3112.2
$ws.Range("J31").Value = 6035
$ws.Range("K31").Value = 3112.2
$ws.Range("L31").Value = 6035
$ws.Range("M31").Value = -2820.2
$ws.Range("N31").Value = -6619
$ws.Range("H37").Value = 3599.3333
$ws.Range("I37").Value = 3112.2
$ws.Range("J37").Value = 6035
$ws.Range("K37").Value = 3112.2
$ws.Range("L37").Value = 6035
$ws.Range("M37").Value = -2835.2
$ws.Range("N37").Value = -6589
$ws.Range("H113").Value = 5087.375
$ws.Range("I113").Value = 2539.8
$ws.Range("J113").Value = 9333.333000000001
$ws.Range("K113").Value = 2539.8
$ws.Range("L113").Value = 9333.333000000001
$ws.Range("M113").Value = -369.8000000000002
$ws.Range("H122").Value = 2153.1667
$ws.Range("I122").Value = 1783.8
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5351.4
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2901.4
$ws.Range("H132").Value = 8999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 8999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 26997
$ws.Range("M132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2600
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 2800
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 2800
$ws.Range("M32").Value = -1683
$ws.Range("N32").Value = -3434
$ws.Range("H68").Value = 6699
$ws.Range("I68").Value = 2747.5
$ws.Range("J68").Value = 9333.333000000001
$ws.Range("K68").Value = 2747.5
$ws.Range("L68").Value = 9333.333000000001
$ws.Range("M68").Value = -1998.5
$ws.Range("N68").Value = -10831.333
$ws.Range("H71").Value = 6699
$ws.Range("I71").Value = 2747.5
$ws.Range("J71").Value = 9333.333000000001
$ws.Range("K71").Value = 13737.5
$ws.Range("L71").Value = 46666.665
$ws.Range("M71").Value = -9993.5
$ws.Range("N71").Value = -54154.665
$ws.Range("H122").Value = 3485.75
$ws.Range("I122").Value = 3485.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10457.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8007.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 12500
$ws.Range("I43").Value = 12500
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 12500
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -12351
$ws.Range("H81").Value = 5555.727
$ws.Range("I81").Value = 4111.3
$ws.Range("J81").Value = 20000
$ws.Range("K81").Value = 8222.6
$ws.Range("L81").Value = 40000
$ws.Range("M81").Value = -7161.6
$ws.Range("H84").Value = 5555.727
$ws.Range("I84").Value = 4111.3
$ws.Range("J84").Value = 20000
$ws.Range("K84").Value = 41113
$ws.Range("L84").Value = 200000
$ws.Range("M84").Value = -35809
$ws.Range("H107").Value = 2999.75
$ws.Range("I107").Value = 2999.75
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 8999.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -7079.25
$ws.Range("H126").Value = 6779.8
$ws.Range("I126").Value = 3932.6667
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 11798.0001
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -9328.000100000001
$ws.Range("H130").Value = 14000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 14000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 14000
$ws.Range("N130").Value = -24040
